$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = '69.713.89'
$ws.Range("E2").Value = '  -0.92%  '

# Row 3: 'Ethereum'
$ws.Range("D3").Value = '3.513.47'
$ws.Range("E3").Value = '  -2.46%  '

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = '  -0.07%  '

# Row 5: 'BNB'
$ws.Range("D5").Value = '''587.66'
$ws.Range("E5").Value = '  +1.34%  '

# Row 6: 'Solana'
$ws.Range("D6").Value = '''183.82'
$ws.Range("E6").Value = '  -3.29%  '

# Row 7: 'LidoStakedEther'
$ws.Range("D7").Value = '3.503.04'
$ws.Range("E7").Value = '  -2.63%  '

# Row 8: 'XRP'
$ws.Range("E8").Value = '  -3.05%  '

# Row 10: 'Dogecoin'
$ws.Range("D10").Value = '''0.198'
$ws.Range("E10").Value = '  +7.29%  '

# Row 11: 'Cardano'
$ws.Range("E11").Value = '  -2.63%  '

# Row 12: 'Avalanche'
$ws.Range("D12").Value = '''54.08'
$ws.Range("E12").Value = '  -3.62%  '

# Row 13: 'ShibaInu'
$ws.Range("E13").Value = '  -2.31%  '

# Row 14: 'Polkadot'
$ws.Range("D14").Value = '''9.46'
$ws.Range("E14").Value = '  -2.35%  '

# Row 15: 'WrappedliquidstakedEther2.0'
$ws.Range("D15").Value = '4.069.65'
$ws.Range("E15").Value = '  -2.71%  '

# Row 16: 'Chainlink'
$ws.Range("D16").Value = '''19.32'
$ws.Range("E16").Value = '  -2.45%  '

# Row 17: 'WrappedBTC'
$ws.Range("D17").Value = '69.683.26'
$ws.Range("E17").Value = '  -0.94%  '

# Row 18: 'WrappedEther'
$ws.Range("D18").Value = '3.501.41'
$ws.Range("E18").Value = '  -2.76%  '

# Row 19: 'Uniswap'
$ws.Range("D19").Value = '''12.33'
$ws.Range("E19").Value = '  -2.50%  '

# Row 20: 'TRON'
$ws.Range("E20").Value = '  -1.45%  '

# Row 21: 'BitcoinCash'
$ws.Range("D21").Value = '''534.97'
$ws.Range("E21").Value = '  +9.05%  '

# Row 22: 'Polygon'
$ws.Range("E22").Value = '  -3.50%  '

# Row 23: 'InternetComputer(DFINITY)'
$ws.Range("D23").Value = '''18.31'
$ws.Range("E23").Value = '  -6.55%  '

# Row 24: 'PancakeSwap'
$ws.Range("D24").Value = '''4.59'
$ws.Range("E24").Value = '  +5.14%  '

# Row 25: 'Toncoin'
$ws.Range("E25").Value = '  -1.00%  '

# Row 26: 'Litecoin'
$ws.Range("D26").Value = '''95.66'
$ws.Range("E26").Value = '  -1.09%  '

# Row 27: 'RenderToken' -> 'ImmutableX'
$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = '''2.98'
$ws.Range("E27").Value = '  -0.50%  '

# Row 28: 'ImmutableX' -> 'RenderToken'
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''11.11'
$ws.Range("E28").Value = '  +0.88%  '

# Row 29: 'Filecoin'
$ws.Range("D29").Value = '''9.10'
$ws.Range("E29").Value = '  -2.95%  '

# Row 30: 'EthereumClassic'
$ws.Range("D30").Value = '''32.21'
$ws.Range("E30").Value = '  -0.17%  '

# Row 31: 'NEARProtocol'
$ws.Range("D31").Value = '''7.30'
$ws.Range("E31").Value = '  -4.21%  '

# Row 32: 'Cosmos'
$ws.Range("E32").Value = '  +1.10%  '

# Row 33: 'OKB'
$ws.Range("D33").Value = '''64.01'
$ws.Range("E33").Value = '  -3.33%  '

# Row 34: 'Hedera'
$ws.Range("E34").Value = '  -3.44%  '

# Row 35: 'Bittensor'
$ws.Range("D35").Value = '''545.56'
$ws.Range("E35").Value = '  -5.31%  '

# Row 36: 'Fetch.AI' -> 'TheGraph'
$ws.Range("B36").Value = 'TheGraph'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D36").Value = '''0.410'
$ws.Range("E36").Value = '  +2.93%  '

# Row 37: 'TheGraph' -> 'Fetch.AI'
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '''3.11'
$ws.Range("E37").Value = '  +5.65%  '

# Row 38: 'InjectiveProtocol'
$ws.Range("D38").Value = '''38.12'
$ws.Range("E38").Value = '  -1.73%  '

# Row 39: 'Dai'
$ws.Range("E39").Value = '  -0.15%  '

# Row 40: 'PEPE'
$ws.Range("D40").Value = '0.0₃0763'
$ws.Range("E40").Value = '  -6.00%  '

# Row 41: 'Kaspa' -> 'Stacks'
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''3.40'
$ws.Range("E41").Value = '  -2.21%  '

# Row 42: 'Stacks' -> 'Kaspa'
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.135'
$ws.Range("E42").Value = '  -2.04%  '

# Row 43: 'Maker'
$ws.Range("D43").Value = '3.362.77'
$ws.Range("E43").Value = '  +4.55%  '

# Row 44: 'dogwifhat'
$ws.Range("D44").Value = '''3.11'
$ws.Range("E44").Value = '  -4.63%  '

# Row 45: 'ApeXProtocol'
$ws.Range("D45").Value = '''3.51'
$ws.Range("E45").Value = '  +2.71%  '

# Row 46: 'ThetaToken'
$ws.Range("E46").Value = '  -2.53%  '

# Row 47: 'VeChain'
$ws.Range("D47").Value = '''0.0440'
$ws.Range("E47").Value = '  -1.39%  '

# Row 48: 'Stellar'
$ws.Range("D48").Value = '''0.135'
$ws.Range("E48").Value = '  -2.97%  '

# Row 49: 'THORChain'
$ws.Range("D49").Value = '''8.91'
$ws.Range("E49").Value = '  -7.21%  '

# Row 50: 'FirstDigitalUSD'
$ws.Range("D50").Value = '''0.998'
$ws.Range("E50").Value = '  -0.12%  '

# Row 51: 'Monero'
$ws.Range("D51").Value = '''137.07'
$ws.Range("E51").Value = '  +1.86%  '
